{"js": "// Add in model standard errors to the \"Descriptive Results\" body paragraphs.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items,text\");\nawait context.sync();\n\n// Locate the three BodyText paragraphs by distinctive leading text instead of\n// hard-coded indices, so the script is resilient to minor structural changes.\nlet statsPara = null;   // \"Snowshoe hare densities decreased...\"\nlet foragingPara = null; // \"Averaged by week, snowshoe hare foraging rate...\"\nlet proteinPara = null;  // \"The average protein composition of fecal sample...\"\n\nfor (const p of paragraphs.items) {\n  const t = p.text;\n  if (t.indexOf(\"Snowshoe hare densities decreased over each winter\") === 0) {\n    statsPara = p;\n  } else if (t.indexOf(\"Averaged by week, snowshoe hare foraging rate\") === 0) {\n    foragingPara = p;\n  } else if (t.indexOf(\"The average protein composition of fecal sample\") === 0) {\n    proteinPara = p;\n  }\n}\n\nif (!statsPara || !foragingPara || !proteinPara) {\n  throw new Error(\"Could not locate expected paragraphs in document body.\");\n}\n\nconst newStatsText =\n  \"Snowshoe hare densities decreased over each winter, because they were always predated, with the greatest overwinter declines occurring from 2016 to 2018, after which hare densities were mostly below the median of 0.37 \\u00b1 0.31 hares per hectare (Figure 1A). Over our study period, hares experienced a 24 \\u00b1 10% chance of being predated, with 2020 showing the highest mortality rate (40%) and 2017 showing the lowest (15%; Figure 1B). The median daily temperature across all years was -15.63 \\u00b1 7.55 C (Figure 1C). Temperatures fluctuated week to week, going above and below the median temperature, except during 2016 when all weeks were warmer than the median (Figure 1C). The median snow depth across all winters was 41.19 \\u00b1 15.97 cm. Snow increased over every winter, reaching a maximum of 78.2 cm in 2021 (Figure 1D). The pattern of accumulation differed between winters; some winters experienced gradual increases (e.g., 2018 and 2019) while others experienced dramatic increases over short periods of time (e.g., 2017 and 2020; Figure 1D). In 2016 and 2019, when temperatures were relatively warmer, snow began to melt and decline in march (Figure 1D). After converting snow depth to food availability, we found there to be a median of 29.46 \\u00b1 8.05kg per hectare of soluble willow twigs available to hares (Figure 1E). After using hare density to calculate food availability on a per capita basis, this equated to 89.24 \\u00b1 78.69 kg of soluble willow per hare (Figure 1F).\";\n\nconst newForagingText =\n  \"Averaged by week, snowshoe hare foraging rate was very flexible (9.42 \\u00b1 1.73), ranging from a minimum of 2.01 to 14.93 hours per day (Figure 2A). Food supplementation reduced female foraging effort by 0.72 hours (p = 0). Within the control sample, male hares foraged 0.42 hours more than females (p = 0.002). As day length increased over winter from January to March, hares decreased their foraging rate by 10.75 \\u00b1 1.19 minutes per hour decrease in night length (p = 0). The median protein composition of fecal sample from food supplemented and control individuals were 11.45 \\u00b1 1.79% and 10.48 \\u00b1 1.71%, respectively. Fecal samples collected in January were higher in protein (11.23%) than those collected in March (10.57%; p = 0.019).\";\n\n// Replace the text of the stats paragraph in place (keeps its style/formatting).\nstatsPara.insertText(newStatsText, Word.InsertLocation.replace);\n\n// Merge the protein paragraph's (updated) content into the foraging paragraph,\n// then delete the now-redundant protein paragraph entirely.\nforagingPara.insertText(newForagingText, Word.InsertLocation.replace);\nproteinPara.delete();\n\nawait context.sync();\n", "ps1": "# Add in model standard errors to the \"Descriptive Results\" body paragraphs.\n$d = $word.ActiveDocument\n\n$statsPara = $null\n$foragingPara = $null\n$proteinPara = $null\n\n$count = $d.Paragraphs.Count\nfor ($i = 1; $i -le $count; $i++) {\n    $p = $d.Paragraphs($i)\n    $t = $p.Range.Text\n    if ($t.StartsWith(\"Snowshoe hare densities decreased over each winter\")) {\n        $statsPara = $p\n    } elseif ($t.StartsWith(\"Averaged by week, snowshoe hare foraging rate\")) {\n        $foragingPara = $p\n    } elseif ($t.StartsWith(\"The average protein composition of fecal sample\")) {\n        $proteinPara = $p\n    }\n}\n\nif (-not $statsPara -or -not $foragingPara -or -not $proteinPara) {\n    throw \"Could not locate expected paragraphs in document body.\"\n}\n\n$newStatsText = \"Snowshoe hare densities decreased over each winter, because they were always predated, with the greatest overwinter declines occurring from 2016 to 2018, after which hare densities were mostly below the median of 0.37 \u00b1 0.31 hares per hectare (Figure 1A). Over our study period, hares experienced a 24 \u00b1 10% chance of being predated, with 2020 showing the highest mortality rate (40%) and 2017 showing the lowest (15%; Figure 1B). The median daily temperature across all years was -15.63 \u00b1 7.55 C (Figure 1C). Temperatures fluctuated week to week, going above and below the median temperature, except during 2016 when all weeks were warmer than the median (Figure 1C). The median snow depth across all winters was 41.19 \u00b1 15.97 cm. Snow increased over every winter, reaching a maximum of 78.2 cm in 2021 (Figure 1D). The pattern of accumulation differed between winters; some winters experienced gradual increases (e.g., 2018 and 2019) while others experienced dramatic increases over short periods of time (e.g., 2017 and 2020; Figure 1D). In 2016 and 2019, when temperatures were relatively warmer, snow began to melt and decline in march (Figure 1D). After converting snow depth to food availability, we found there to be a median of 29.46 \u00b1 8.05kg per hectare of soluble willow twigs available to hares (Figure 1E). After using hare density to calculate food availability on a per capita basis, this equated to 89.24 \u00b1 78.69 kg of soluble willow per hare (Figure 1F).\"\n\n$newForagingText = \"Averaged by week, snowshoe hare foraging rate was very flexible (9.42 \u00b1 1.73), ranging from a minimum of 2.01 to 14.93 hours per day (Figure 2A). Food supplementation reduced female foraging effort by 0.72 hours (p = 0). Within the control sample, male hares foraged 0.42 hours more than females (p = 0.002). As day length increased over winter from January to March, hares decreased their foraging rate by 10.75 \u00b1 1.19 minutes per hour decrease in night length (p = 0). The median protein composition of fecal sample from food supplemented and control individuals were 11.45 \u00b1 1.79% and 10.48 \u00b1 1.71%, respectively. Fecal samples collected in January were higher in protein (11.23%) than those collected in March (10.57%; p = 0.019).\"\n\n# Replace the text of the stats paragraph in place (keeps its style/formatting).\n$statsPara.Range.Text = $newStatsText\n\n# Merge the protein paragraph's (updated) content into the foraging paragraph,\n# then delete the now-redundant protein paragraph entirely.\n$foragingPara.Range.Text = $newForagingText\n$proteinPara.Range.Delete()\n"}
